$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells holding plain percentage-like text (e.g. "70%") need an explicit
# text number format, otherwise Excel auto-converts the string into a
# numeric percentage value instead of keeping the literal text.
$ws.Range("E2").Value = "2026-02-09 22:48:25"
$ws.Range("I2").Value = "6.5 mm"
$ws.Range("E3").Value = "2026-02-09 22:48:27"
$ws.Range("I3").Value = "5.2 mm"
$ws.Range("E4").Value = "2026-02-09 22:48:29"
$ws.Range("J4").Value = "1006.9 hPa"
$ws.Range("K4").Value = "11.2 MJ/m2"
$ws.Range("E5").Value = "2026-02-09 22:48:31"
$ws.Range("G5").Value = "122 cm"
$ws.Range("I5").Value = "2.9 mm"
$ws.Range("O5").Value = "-2.5 °C"
$ws.Range("E6").Value = "2026-02-09 22:48:33"
$ws.Range("J6").Value = "1006.8 hPa"
$ws.Range("E7").Value = "2026-02-09 22:48:36"
$ws.Range("I7").Value = "1.1 mm"
$ws.Range("E8").Value = "2026-02-09 22:48:38"
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = "70%"
$ws.Range("I8").Value = "1.3 mm"
$ws.Range("E9").Value = "2026-02-09 22:48:40"
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = "83%"
$ws.Range("E10").Value = "2026-02-09 22:48:43"
$ws.Range("E11").Value = "2026-02-09 22:48:45"
$ws.Range("I11").Value = "1.1 mm"
$ws.Range("O11").Value = "5.2 °C"
$ws.Range("E12").Value = "2026-02-09 22:48:47"
$ws.Range("E13").Value = "2026-02-09 22:48:49"
$ws.Range("I13").Value = "1.4 mm"
$ws.Range("E14").Value = "2026-02-09 22:48:52"
$ws.Range("I14").Value = "1.2 mm"
$ws.Range("E15").Value = "2026-02-09 22:48:54"
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = "82%"
$ws.Range("E16").Value = "2026-02-09 22:48:56"
$ws.Range("I16").Value = "3.6 mm"
$ws.Range("O16").Value = "-3.3 °C"
$ws.Range("E17").Value = "2026-02-09 22:48:59"
$ws.Range("E18").Value = "2026-02-09 22:49:01"
$ws.Range("H18").NumberFormat = "@"
$ws.Range("H18").Value = "83%"
$ws.Range("E19").Value = "2026-02-09 22:49:04"
$ws.Range("I19").Value = "0.5 mm"
$ws.Range("E20").Value = "2026-02-09 22:49:06"
$ws.Range("I20").Value = "1.6 mm"
$ws.Range("M20").Value = "-1.5 °C 22:29 TU"
$ws.Range("O20").Value = "-4.1 °C"
$ws.Range("E21").Value = "2026-02-09 22:49:08"
$ws.Range("I21").Value = "1.4 mm"
$ws.Range("E22").Value = "2026-02-09 22:49:11"
$ws.Range("G22").Value = "121 cm"
$ws.Range("L22").Value = "33.5 km/h - 305º 22:15 TU"
$ws.Range("E23").Value = "2026-02-09 22:49:13"
$ws.Range("G23").Value = "173 cm"
$ws.Range("I23").Value = "5.2 mm"
$ws.Range("E24").Value = "2026-02-09 22:49:16"
$ws.Range("J24").Value = "1008.4 hPa"
$ws.Range("E25").Value = "2026-02-09 22:49:18"
$ws.Range("H25").NumberFormat = "@"
$ws.Range("H25").Value = "76%"
$ws.Range("I25").Value = "2.2 mm"
$ws.Range("E26").Value = "2026-02-09 22:49:20"
$ws.Range("H26").NumberFormat = "@"
$ws.Range("H26").Value = "83%"
$ws.Range("J26").Value = "1006.6 hPa"
$ws.Range("O26").Value = "2.7 °C"
$ws.Range("E27").Value = "2026-02-09 22:49:23"
$ws.Range("G27").Value = "172 cm"
$ws.Range("H27").NumberFormat = "@"
$ws.Range("H27").Value = "84%"
$ws.Range("I27").Value = "2.5 mm"
$ws.Range("E28").Value = "2026-02-09 22:49:25"
$ws.Range("J28").Value = "1006.9 hPa"
$ws.Range("E29").Value = "2026-02-09 22:49:28"
$ws.Range("E30").Value = "2026-02-09 22:49:30"
$ws.Range("H30").NumberFormat = "@"
$ws.Range("H30").Value = "87%"
$ws.Range("E31").Value = "2026-02-09 22:49:32"
$ws.Range("J31").Value = "1006.2 hPa"
$ws.Range("E32").Value = "2026-02-09 22:49:35"
$ws.Range("I32").Value = "2.1 mm"
$ws.Range("E33").Value = "2026-02-09 22:49:37"
$ws.Range("I33").Value = "1.4 mm"
$ws.Range("E34").Value = "2026-02-09 22:49:40"
$ws.Range("I34").Value = "1.1 mm"
$ws.Range("E35").Value = "2026-02-09 22:49:42"
$ws.Range("J35").Value = "1008.8 hPa"
$ws.Range("M35").Value = "9.4 °C 22:29 TU"
$ws.Range("O35").Value = "5.6 °C"
$ws.Range("E36").Value = "2026-02-09 22:49:44"
$ws.Range("J36").Value = "1007.0 hPa"
$ws.Range("E37").Value = "2026-02-09 22:49:47"
$ws.Range("I37").Value = "0.2 mm"
$ws.Range("E38").Value = "2026-02-09 22:49:49"
$ws.Range("E39").Value = "2026-02-09 22:49:51"
$ws.Range("O39").Value = "-3.1 °C"
$ws.Range("E40").Value = "2026-02-09 22:49:54"
$ws.Range("I40").Value = "2.1 mm"
$ws.Range("E41").Value = "2026-02-09 22:49:56"
$ws.Range("H41").NumberFormat = "@"
$ws.Range("H41").Value = "59%"
$ws.Range("E42").Value = "2026-02-09 22:49:58"
$ws.Range("E43").Value = "2026-02-09 22:50:01"
$ws.Range("H43").NumberFormat = "@"
$ws.Range("H43").Value = "74%"
$ws.Range("I43").Value = "0.3 mm"
$ws.Range("E44").Value = "2026-02-09 22:50:03"
$ws.Range("I44").Value = "2.8 mm"
$ws.Range("E45").Value = "2026-02-09 22:50:05"
$ws.Range("I45").Value = "2.5 mm"
$ws.Range("E46").Value = "2026-02-09 22:50:08"
$ws.Range("H46").NumberFormat = "@"
$ws.Range("H46").Value = "76%"
$ws.Range("J46").Value = "1008.6 hPa"
